# Weekly update: insert a new daily price record for Plátano (Agrícola del
# Norte S.A. de Arica) before the existing row 285, pushing the remainder of
# the table (old rows 285-384) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 285 - this shifts rows 285:384 down to
# 286:385 and extends the used range to A1:T385, matching the existing
# row's formatting (including the date cell's style) along the way.
$ws.Rows("285:285").Insert()

# Populate the new row with the inserted record.
$ws.Range("A285").Value = 1
$ws.Range("B285").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C285").Value = "Arica y Parinacota"
$ws.Range("D285").Value = 45027
$ws.Range("E285").Value = 15
$ws.Range("F285").Value = "Fruta"
$ws.Range("G285").Value = 100108
$ws.Range("H285").Value = "Tropicales y subtropicales"
$ws.Range("I285").Value = 100108006
$ws.Range("J285").Value = "Plátano"
$ws.Range("K285").Value = "Sin especificar"
$ws.Range("L285").Value = "Pintón"
$ws.Range("M285").Value = 120
$ws.Range("N285").Value = 19000
$ws.Range("O285").Value = 20000
$ws.Range("P285").Value = 19500
$ws.Range("Q285").Value = "$/caja 20 kilos"
$ws.Range("R285").Value = "Ecuador"
$ws.Range("S285").Value = 975
$ws.Range("T285").Value = 20
